$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update regression result values (rows 2-23) with the latest data
$ws.Range("A2").Value = "0.568***"
$ws.Range("A3").Value = "(0.004)"
$ws.Range("A4").Value = "-0.085***"
$ws.Range("A5").Value = "(0.002)"
$ws.Range("A6").Value = "-0.099***"
$ws.Range("A7").Value = "(0.024)"
$ws.Range("A8").Value = "-0.001***"
$ws.Range("A9").Value = "(0.000)"
# Prefix with an apostrophe so the numeric-looking string "-0.001" is kept
# as text (matching the rest of the column) instead of being parsed as a
# number; reset the style afterwards so no stray quote-prefix format is left.
$ws.Range("A10").Value = "'-0.001"
$ws.Cells.Item(10, 1).Style = "Normal"
$ws.Range("A11").Value = "(0.000)"
$ws.Range("A12").Value = "-0.000***"
$ws.Range("A13").Value = "(0.000)"
$ws.Range("A14").Value = "-0.002***"
$ws.Range("A15").Value = "(0.000)"
$ws.Range("A16").Value = "-0.001***"
$ws.Range("A17").Value = "(0.000)"
$ws.Range("A18").Value = "-0.000***"
$ws.Range("A19").Value = "(0.000)"
$ws.Range("A20").Value = "-0.014***"
$ws.Range("A21").Value = "(0.000)"
$ws.Range("A22").Value = "-0.000***"
$ws.Range("A23").Value = "(0.000)"

# Remove the now-obsolete trailing rows (24-31) from the old, longer data set
$ws.Range("A24:A31").EntireRow.Delete()
